$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J header block (row 1) - copy formatting of column I row 1 (plain top style)
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# Row 2 header label "SE-43" with same style as the other week headers
$ws.Range("I2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$ws.Range("J2").Value = "SE-43"

# Row 3 sub-header "24oct-30oct" with same style as the other date-range cells
$ws.Range("I3").Copy()
$ws.Range("J3").PasteSpecial(-4122)
$ws.Range("J3").Value = "24oct-30oct"

# Rows 4-8, 11-12, 14-15: empty data cells, styled like the rest of the data area
$ws.Range("D4").Copy()
$ws.Range("J4:J8").PasteSpecial(-4122)
$ws.Range("J11:J12").PasteSpecial(-4122)
$ws.Range("J14:J15").PasteSpecial(-4122)

# Rows 9, 10, 13: new counts for SE-43 (Chumbivilcas, Cusco, La Convención)
$ws.Range("J9").Value = 1
$ws.Range("J10").Value = 2
$ws.Range("J13").Value = 1

# Row 16 totals row - copy formatting of the adjacent total-row style cell
$ws.Range("I16").Copy()
$ws.Range("J16").PasteSpecial(-4122)

# Update view: select A1:J17, leave the active cell on J17 (bottom-right,
# the newly filled corner), and zoom to 78%
$ws.Application.CutCopyMode = $false
$ws.Range("A1:J17").Select()
$ws.Range("J17").Activate()
$excel.ActiveWindow.Zoom = 78
